# PC3.xlsx "fix data and script" edit:
#  - sheet grows from A1:T3 to A1:X4
#  - T1 header "CC50" -> "CC50.x"; new headers U1:X1 = CC50.y / Lower / Upper / SE
#  - row 2 ("OA") replaced by new "DMSO" stats
#  - row 3 ("DMSO") replaced by new "MMAE" stats (+ new U3:X3 values)
#  - row 4 is a brand-new "MMP62" stats row (A4:X4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers (only T1 changes text; U1:X1 are new columns)
$headers = @(
    "Drug","F val","p-val","Slope","LL","UL","ED50","Slope SE","LL SE","UL SE",
    "ED50 SE","Slope t-val","LL t-val","UL t-val","ED50 t-val","Slope p-val",
    "LL p-val","UL p-val","ED50 p-val","CC50.x","CC50.y","Lower","Upper","SE"
)
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Rows 2-4 - data (columns A..X == indices 1..24). $null marks a cell that
# should be left blank (no value written) for that row.
$rows = @(
    @("DMSO", 0.546, 0.705, -3.491, 71.675, 100.052, 0.151, 9.481999999999999,
      12.97, 4.926, 0.08500000000000001, -0.368, 5.526, 20.309, 1.783, 0.717,
      0, 0, 0.091, 100, $null, $null, $null, $null),
    @("MMAE", 1.191, 0.355, 1.194, 28.102, 63.048, 0.001, 1.186, 4.984, 4.354,
      0.001, 1.007, 5.639, 14.48, 1.262, 0.327, 0, 0, 0.222,
      [double]"2.753623717932641e-012", 0.001661953671243145,
      -0.002253639655625769, 0.005577546998112059, 0.0014102902015021),
    @("MMP62", 1.342, 0.3, 1.441, 24.495, 78.867, 0.332, 0.481, 3.613, 3.566,
      0.117, 2.996, 6.78, 22.114, 2.84, 0.007, 0, 0, 0.01,
      0.3600000000027537, 0.4719423107703062, 0.3025823117782635,
      0.6413023097623489, 0.06099886458226141)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowVals = $rows[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        if ($null -ne $rowVals[$c]) {
            $ws.Cells.Item($r + 2, $c + 1).Value = $rowVals[$c]
        }
    }
}
